# Apply the "gh-pages output generated at 456a3b4" update:
#  - bump a handful of "想去人数" (interest-count) numbers in the
#    展览 (sheet index 1), 演出 (sheet index 2) and 全部类型 (sheet index 4)
#    worksheets
#  - insert one brand-new event row (2024-08-03 ACG concert) into both the
#    演出 sheet (before its current last row) and the 全部类型 sheet
#    (before the matching row), shifting the row that used to be last down
#    by one.

$wb = $excel.ActiveWorkbook

function Set-CountCell {
    param(
        $ws,
        [int]$row,
        [double]$value
    )
    $ws.Range("F" + $row).Value = $value
}

function Insert-EventRow {
    param(
        $ws,
        [int]$rowIndex,      # row number the new row will occupy
        [int]$indexValue,    # value for column A (sequential index)
        [string]$startDate,  # column B - must stay text, not become a date serial
        [string]$title,      # column C
        [string]$place,      # column D
        [string]$timeRange,  # column E
        [double]$wantCount,  # column F
        [double]$minPrice,   # column G
        [string]$link,       # column H
        [string]$cover       # column I
    )

    # Shift everything from $rowIndex down by one row, copying the format
    # of the row immediately above so the newly created row looks the same
    # as its neighbours.
    $ws.Rows.Item($rowIndex).Insert()

    $donorRow = $rowIndex - 1

    # Column A needs the bold/centered/bordered "index" style used by every
    # other cell in column A.
    $ws.Range("A" + $donorRow).Copy() | Out-Null
    $ws.Range("A" + $rowIndex).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = $false
    $ws.Range("A" + $rowIndex).Value = $indexValue

    # Column B holds a date-shaped string ("2024-08-03"); Excel would
    # silently convert a bare string like that into a date serial number,
    # so force text entry and then strip the formatting mark back off so
    # the cell ends up plain/unstyled just like its neighbours.
    $ws.Range("B" + $rowIndex).NumberFormat = "@"
    $ws.Range("B" + $rowIndex).Value = $startDate
    $ws.Range("C" + $donorRow).Copy() | Out-Null
    $ws.Range("B" + $rowIndex).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $excel.CutCopyMode = $false

    $ws.Range("C" + $rowIndex).Value = $title
    $ws.Range("D" + $rowIndex).Value = $place
    $ws.Range("E" + $rowIndex).Value = $timeRange
    $ws.Range("F" + $rowIndex).Value = $wantCount
    $ws.Range("G" + $rowIndex).Value = $minPrice
    $ws.Range("H" + $rowIndex).Value = $link
    $ws.Range("I" + $rowIndex).Value = $cover
}

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions) - just the want-count bumps
# ---------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
Set-CountCell $wsExhibit 4  1755
Set-CountCell $wsExhibit 5  443
Set-CountCell $wsExhibit 9  307
Set-CountCell $wsExhibit 10 1698
Set-CountCell $wsExhibit 11 340
Set-CountCell $wsExhibit 12 1399
Set-CountCell $wsExhibit 13 792
Set-CountCell $wsExhibit 14 321
Set-CountCell $wsExhibit 15 663
Set-CountCell $wsExhibit 16 12650
Set-CountCell $wsExhibit 17 12672
Set-CountCell $wsExhibit 18 937
Set-CountCell $wsExhibit 21 299
Set-CountCell $wsExhibit 23 509
Set-CountCell $wsExhibit 24 1981
Set-CountCell $wsExhibit 25 23
Set-CountCell $wsExhibit 28 663

# ---------------------------------------------------------------------
# Sheet "演出" (performances) - want-count bump + new row insertion
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
Set-CountCell $wsShow 5 73

Insert-EventRow $wsShow 11 10 "2024-08-03" `
    "广州·【暑期5折】《忱宴·渐渐被你吸引》热血动漫二次元ACG演唱会" `
    "东风中路299号 广州中山纪念堂" `
    "2024.08.03 20:00-08.03 21:40" `
    2 50 `
    "https://show.bilibili.com/platform/detail.html?id=85917" `
    "//i1.hdslb.com/bfs/openplatform/202405/won43hte1715675570347.jpeg"

# the row that used to be row 11 is now row 12 - renumber its index cell
$wsShow.Range("A12").Value = 11

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types) - want-count bumps + new row insertion
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
Set-CountCell $wsAll 6  1755
Set-CountCell $wsAll 7  443
Set-CountCell $wsAll 14 307
Set-CountCell $wsAll 15 1698
Set-CountCell $wsAll 16 340
Set-CountCell $wsAll 17 1399
Set-CountCell $wsAll 18 792
Set-CountCell $wsAll 19 321
Set-CountCell $wsAll 20 73
Set-CountCell $wsAll 21 663
Set-CountCell $wsAll 22 12650
Set-CountCell $wsAll 23 12672
Set-CountCell $wsAll 24 937
Set-CountCell $wsAll 27 299
Set-CountCell $wsAll 29 509
Set-CountCell $wsAll 32 1981
Set-CountCell $wsAll 33 23
Set-CountCell $wsAll 38 663

Insert-EventRow $wsAll 40 39 "2024-08-03" `
    "广州·【暑期5折】《忱宴·渐渐被你吸引》热血动漫二次元ACG演唱会" `
    "东风中路299号 广州中山纪念堂" `
    "2024.08.03 20:00-08.03 21:40" `
    2 50 `
    "https://show.bilibili.com/platform/detail.html?id=85917" `
    "//i1.hdslb.com/bfs/openplatform/202405/won43hte1715675570347.jpeg"

# the row that used to be row 40 is now row 41 - renumber its index cell
$wsAll.Range("A41").Value = 40
